$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data A:E -> B:F
$ws.Columns("A:A").Insert()

# Copy formatting from the (now shifted) header cell B1 onto the new A1 header cell
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set header text for new ID column
$ws.Range("A1").Value = "ID"

# Fill in the row identifiers for the new ID column
$ids = @("Hb 2","Hb 3","S 24","S 28","Hb 107","Hb 66","Hb 69","Hb 95","Hb 99","Hb 92","Hb 40","Hb 41","S 11","Hb 57","S 21","S 22","S 3","S 4","S 5","Hb 74","Hb 79","Hb 32","S 15","S 16")
for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $ids[$i]
}
